$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update inserts two new price records (Castle Brite Primera /
# Segunda, $/caja 15 kilos) at the top of the data block, pushing the
# existing rows 7-20 down to rows 9-22 and extending the used range to
# A1:T22.
$ws.Rows("7:8").Insert()

# Row 7 - new record
$ws.Range("A7").Value2 = 5
$ws.Range("B7").Value2 = "Macroferia Regional de Talca"
$ws.Range("C7").Value2 = "Maule"
$ws.Range("D7").Value2 = 44533
$ws.Range("E7").Value2 = 7
$ws.Range("F7").Value2 = "Fruta"
$ws.Range("G7").Value2 = 100103
$ws.Range("H7").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I7").Value2 = 100103003
$ws.Range("J7").Value2 = "Damasco"
$ws.Range("K7").Value2 = "Castle Brite"
$ws.Range("L7").Value2 = "Primera"
$ws.Range("M7").Value2 = 170
$ws.Range("N7").Value2 = 18000
$ws.Range("O7").Value2 = 18000
$ws.Range("P7").Value2 = 18000
$ws.Range("Q7").Value2 = "`$/caja 15 kilos"
$ws.Range("R7").Value2 = "Región de O'Higgins"
$ws.Range("S7").Value2 = 1200
$ws.Range("T7").Value2 = 15

# Row 8 - new record
$ws.Range("A8").Value2 = 5
$ws.Range("B8").Value2 = "Macroferia Regional de Talca"
$ws.Range("C8").Value2 = "Maule"
$ws.Range("D8").Value2 = 44533
$ws.Range("E8").Value2 = 7
$ws.Range("F8").Value2 = "Fruta"
$ws.Range("G8").Value2 = 100103
$ws.Range("H8").Value2 = "Frutos de hueso (carozo)"
$ws.Range("I8").Value2 = 100103003
$ws.Range("J8").Value2 = "Damasco"
$ws.Range("K8").Value2 = "Castle Brite"
$ws.Range("L8").Value2 = "Segunda"
$ws.Range("M8").Value2 = 100
$ws.Range("N8").Value2 = 14000
$ws.Range("O8").Value2 = 14000
$ws.Range("P8").Value2 = 14000
$ws.Range("Q8").Value2 = "`$/caja 15 kilos"
$ws.Range("R8").Value2 = "Región de O'Higgins"
$ws.Range("S8").Value2 = 933
$ws.Range("T8").Value2 = 15
